# Applies the commit "Connection zu, Benutzerverwaltungssystem muss
# abgelegt werden, erwähnen in Systemarchitektur" to the
# Systemarchitektur documentation:
#   1. "Opennes?" -> "Open-Nes?"
#   2. Reworks the design-tool paragraph: describes how Modeler/Workbench
#      generate the ER model & scripts (instead of the old sentence about
#      the logical ER design being skipped), and appends the
#      "beziehungsweise Nachbearbeitungen" qualifier.
#   3. Fixes a stray double space ("Anbindung  und" -> "Anbindung und").
#   4. Adds two new bulleted paragraphs documenting the BAC1 interfaces
#      (Authorisierung / Importmodule), right after the
#      "Datenbankdesigntool" bullet.

$d = $word.ActiveDocument

# 1) "Opennes?" -> "Open-Nes?"
$ok1 = $d.Content.Find.Execute(
    "Opennes",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Open-Nes", 2)
Write-Output "step1 Opennes->Open-Nes: $ok1"

# 2a) Rewrite the "graphischen Tabellenentwurf" sentence describing the
#     design-tool comparison, dropping the ER-design paragraph in favour
#     of a description of how both tools generate the ER model / scripts.
$ok2 = $d.Content.Find.Execute(
    "Für den graphischen Tabellenentwurf bieten beide Tools Unterstützung. Auf den Vorteil des logischen (ER) Designs, welches der Oracle Modeler unterstützt wird bewusst verzichtet, da der Aufwand die Überführung in ein Tabellenmodell zu machen geringer ist, als die Skriptdateien zur Anlage der Tabellen, Indizes, Einschränkungen bezüglich referentieller Integrität und die automatische Vergabe eindeutiger Schlüssel nachzubearbeiten. ",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Für den graphischen Entwurf des Entity-Relation-Modells bieten beide Tools Unterstützung. Modeler wie Workbench generieren aus dem ER-Modell sowohl die graphische Übersicht der Tabellen und Schlüssel, als auch Generation der Skriptdateien zur Anlage der Tabellen, Indizes, Einschränkungen bezüglich referentieller Integrität und die automatische Vergabe eindeutiger Schlüssel. ",
    2)
Write-Output "step2a ER-model sentence rewrite: $ok2"

# 2b) Append the "beziehungsweise Nachbearbeitungen" qualifier.
$ok3 = $d.Content.Find.Execute(
    "einzusetzen und damit Kompatibilitätsprobleme zu vermeiden.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "einzusetzen und damit Kompatibilitätsprobleme beziehungsweise Nachbearbeitungen zu vermeiden.",
    2)
Write-Output "step2b Nachbearbeitungen qualifier: $ok3"

# 3) Fix the double space in "Anbindung  und" -> "Anbindung und"
$ok4 = $d.Content.Find.Execute(
    "Anbindung  und",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Anbindung und",
    2)
Write-Output "step3 double-space fix: $ok4"

# 4) Insert the two new "Schnittstellen zu BAC1" bullet paragraphs right
#    after the "Datenbankdesigntool" bullet and before the
#    "Die Anbindung und die Veröffentlichung" bullet.

# Locate the "Datenbankdesigntool" paragraph.
$dbDesignPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Datenbankdesigntool:*") {
        $dbDesignPara = $p
        break
    }
}
Write-Output "step4 found Datenbankdesigntool paragraph: $($dbDesignPara -ne $null)"

$authText = "Schnittstellen zu BAC1 (Authorisierung): der bisherige Entwurf des API ist in einer generischen Form implementiert, sodass einerseits ein Austausch des Authorisierungsmoduls möglich ist (zum Beispiel in einem ersten Schritt realisiert als Benutzername/Passwort Zugang). Eine Abhängigkeit zum Datenmodell wurde identifiziert: die Form der Benutzerkennung des externen Rollensystems muss in den Tabellen untergebracht werden könne (Details siehe Abschnitt Datenmodell). "
$importText = "Schnittstellen zu BAC1 (Importmodule): gemeinsame Basis für die API und Importmodule ist das Datenmodell (Details zu notwendigen Erweiterungen siehe im Abschnitt Datenmodell), programmtechnisch wurden keine Überschneidungen identifiziert. "

if ($dbDesignPara -ne $null) {
    # New list paragraph #1: Authorisierung.
    $dbDesignPara.Range.InsertParagraphAfter()
    $authPara = $dbDesignPara.Next()
    $authPara.Range.Text = $authText

    # New list paragraph #2: Importmodule.
    $authPara.Range.InsertParagraphAfter()
    $importPara = $authPara.Next()
    $importPara.Range.Text = $importText

    Write-Output "step4 inserted paragraphs: [$($authPara.Range.Text)] / [$($importPara.Range.Text)]"
}
